$wb = $excel.ActiveWorkbook

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
# xlContinuous=1, xlThin=2, xlLineStyleNone=-4142
# xlPasteFormats=-4122

function Set-TopBottomBorder($rng) {
    $rng.ClearFormats()
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(7).LineStyle = -4142
}

function Set-TopRightBottomBorder($rng) {
    $rng.ClearFormats()
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(10).LineStyle = 1
    $rng.Borders.Item(10).Weight = 2
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(7).LineStyle = -4142
}

# --- Sheet "quality_comparison": build the two new border styles once ---
$ws1 = $wb.Worksheets.Item("quality_comparison")
Set-TopBottomBorder $ws1.Range("C1")
Set-TopRightBottomBorder $ws1.Range("D1")
$ws1.Range("C2").Value = "approach"

# --- Sheet "computational_comparison": reuse the same styles via format copy,
#     so no duplicate cellXfs/border records get created ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

$ws1.Range("C1").Copy() | Out-Null
$ws2.Range("C1").PasteSpecial(-4122) | Out-Null
$ws1.Range("D1").Copy() | Out-Null
$ws2.Range("D1").PasteSpecial(-4122) | Out-Null
$ws1.Range("C1").Copy() | Out-Null
$ws2.Range("F1").PasteSpecial(-4122) | Out-Null
$ws1.Range("D1").Copy() | Out-Null
$ws2.Range("G1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"
$ws2.Range("G5").ClearContents()

Write-Output "done"
